$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 85 (shifts existing rows 85+ down by 3)
$ws.Range("A85:T87").EntireRow.Insert()

# Fill the new rows with weekly data
# Row 85
$ws.Cells.Item(85, 1).Value = 2
$ws.Cells.Item(85, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44846
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100107
$ws.Cells.Item(85, 8).Value = "Otros"
$ws.Cells.Item(85, 9).Value = 100107002
$ws.Cells.Item(85, 10).Value = "Chirimoya"
$ws.Cells.Item(85, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(85, 12).Value = "Especial"
$ws.Cells.Item(85, 13).Value = 400
$ws.Cells.Item(85, 14).Value = 19000
$ws.Cells.Item(85, 15).Value = 20000
$ws.Cells.Item(85, 16).Value = 19500
$ws.Cells.Item(85, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(85, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(85, 19).Value = 1950
$ws.Cells.Item(85, 20).Value = 10

# Row 86
$ws.Cells.Item(86, 1).Value = 2
$ws.Cells.Item(86, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = 44846
$ws.Cells.Item(86, 5).Value = 4
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100107
$ws.Cells.Item(86, 8).Value = "Otros"
$ws.Cells.Item(86, 9).Value = 100107002
$ws.Cells.Item(86, 10).Value = "Chirimoya"
$ws.Cells.Item(86, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 400
$ws.Cells.Item(86, 14).Value = 17000
$ws.Cells.Item(86, 15).Value = 18000
$ws.Cells.Item(86, 16).Value = 17500
$ws.Cells.Item(86, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(86, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(86, 19).Value = 1750
$ws.Cells.Item(86, 20).Value = 10

# Row 87
$ws.Cells.Item(87, 1).Value = 2
$ws.Cells.Item(87, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44846
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = "Fruta"
$ws.Cells.Item(87, 7).Value = 100107
$ws.Cells.Item(87, 8).Value = "Otros"
$ws.Cells.Item(87, 9).Value = 100107002
$ws.Cells.Item(87, 10).Value = "Chirimoya"
$ws.Cells.Item(87, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(87, 12).Value = "Segunda"
$ws.Cells.Item(87, 13).Value = 300
$ws.Cells.Item(87, 14).Value = 14000
$ws.Cells.Item(87, 15).Value = 15000
$ws.Cells.Item(87, 16).Value = 14500
$ws.Cells.Item(87, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(87, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(87, 19).Value = 1450
$ws.Cells.Item(87, 20).Value = 10
